$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '26.808.45'
Set-TextValue $ws.Range('E2') '  +4.21%  '
Set-TextValue $ws.Range('D3') '1.873.92'
Set-TextValue $ws.Range('E3') '  +3.13%  '
Set-TextValue $ws.Range('E4') '  -0.01%  '
Set-TextValue $ws.Range('D5') '277.47'
Set-TextValue $ws.Range('E5') '  -0.01%  '
Set-TextValue $ws.Range('E6') '  +0.00%  '
Set-TextValue $ws.Range('D7') '0.5291'
Set-TextValue $ws.Range('E7') '  +3.90%  '
Set-TextValue $ws.Range('D8') '0.3421'
Set-TextValue $ws.Range('E8') '  -2.82%  '
Set-TextValue $ws.Range('E9') '  +4.08%  '
Set-TextValue $ws.Range('D11') '0.8049'
Set-TextValue $ws.Range('E11') '  -2.52%  '
Set-TextValue $ws.Range('D12') '0.07756'
Set-TextValue $ws.Range('E12') '  -1.83%  '
Set-TextValue $ws.Range('D13') '1.867.22'
Set-TextValue $ws.Range('E13') '  +2.72%  '
Set-TextValue $ws.Range('B14') 'Litecoin'
Set-TextValue $ws.Range('C14') 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range('D14') '90.46'
Set-TextValue $ws.Range('E14') '  +3.20%  '
Set-TextValue $ws.Range('B15') 'Polkadot'
Set-TextValue $ws.Range('C15') 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D15') '5.187'
Set-TextValue $ws.Range('E15') '  +2.28%  '
Set-TextValue $ws.Range('D16') '14.58'
Set-TextValue $ws.Range('E16') '  +3.54%  '
Set-TextValue $ws.Range('D17') '0.9998'
Set-TextValue $ws.Range('E17') '  -0.07%  '
Set-TextValue $ws.Range('D18') '0.000008056'
Set-TextValue $ws.Range('E18') '  +0.30%  '
Set-TextValue $ws.Range('D19') '1.000'
Set-TextValue $ws.Range('E19') '  -0.03%  '
Set-TextValue $ws.Range('D20') '26.823.41'
Set-TextValue $ws.Range('E20') '  +4.09%  '
Set-TextValue $ws.Range('D21') '2.100.02'
Set-TextValue $ws.Range('E21') '  +1.56%  '
Set-TextValue $ws.Range('D22') '4.750'
Set-TextValue $ws.Range('E22') '  +0.20%  '
Set-TextValue $ws.Range('D24') '6.173'
Set-TextValue $ws.Range('E24') '  +1.12%  '
Set-TextValue $ws.Range('D25') '2.380'
Set-TextValue $ws.Range('E25') '  +7.57%  '
Set-TextValue $ws.Range('D26') '146.62'
Set-TextValue $ws.Range('E26') '  +3.35%  '
Set-TextValue $ws.Range('B27') 'EthereumClassic'
Set-TextValue $ws.Range('C27') 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D27') '17.35'
Set-TextValue $ws.Range('E27') '  +1.41%  '
Set-TextValue $ws.Range('B28') 'Toncoin'
Set-TextValue $ws.Range('C28') 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D28') '1.662'
Set-TextValue $ws.Range('E28') '  -0.62%  '
Set-TextValue $ws.Range('D29') '113.47'
Set-TextValue $ws.Range('E29') '  +3.80%  '
Set-TextValue $ws.Range('D30') '4.348'
Set-TextValue $ws.Range('E30') '  +0.42%  '
Set-TextValue $ws.Range('D31') '4.318'
Set-TextValue $ws.Range('E31') '  +2.05%  '
Set-TextValue $ws.Range('D32') '0.08918'
Set-TextValue $ws.Range('E32') '  +1.68%  '
Set-TextValue $ws.Range('D33') '0.04933'
Set-TextValue $ws.Range('E33') '  +0.93%  '
Set-TextValue $ws.Range('D35') '0.7311'
Set-TextValue $ws.Range('E35') '  +0.61%  '
Set-TextValue $ws.Range('D36') '2.884'
Set-TextValue $ws.Range('E36') '  +0.53%  '
Set-TextValue $ws.Range('E37') '  +4.92%  '
Set-TextValue $ws.Range('B38') 'RenderToken'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D38') '2.330'
Set-TextValue $ws.Range('E38') '  -1.80%  '
Set-TextValue $ws.Range('B39') 'VeChain'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D39') '0.01858'
Set-TextValue $ws.Range('E39') '  +0.24%  '
Set-TextValue $ws.Range('D40') '0.5159'
Set-TextValue $ws.Range('E40') '  -0.37%  '
Set-TextValue $ws.Range('D41') '0.9521'
Set-TextValue $ws.Range('E41') '  -1.47%  '
Set-TextValue $ws.Range('D42') '116.08'
Set-TextValue $ws.Range('E42') '  +4.71%  '
Set-TextValue $ws.Range('D43') '6.165'
Set-TextValue $ws.Range('E43') '  -0.95%  '
Set-TextValue $ws.Range('D44') '8.128'
Set-TextValue $ws.Range('E44') '  +1.49%  '
Set-TextValue $ws.Range('D45') '0.9997'
Set-TextValue $ws.Range('D46') '0.4478'
Set-TextValue $ws.Range('E46') '  -1.93%  '
Set-TextValue $ws.Range('D47') '0.1342'
Set-TextValue $ws.Range('E47') '  -1.64%  '
Set-TextValue $ws.Range('D48') '9.352'
Set-TextValue $ws.Range('E48') '  +1.68%  '
Set-TextValue $ws.Range('D49') '36.39'
Set-TextValue $ws.Range('E49') '  -0.14%  '
Set-TextValue $ws.Range('D50') '0.05940'
Set-TextValue $ws.Range('E50') '  +1.68%  '
Set-TextValue $ws.Range('D51') '1.494'
Set-TextValue $ws.Range('E51') '  -0.47%  '
